$d = $word.ActiveDocument

# Both affected spots read "MSG " (the literal message-code marker, followed
# by a space) in the document's visible text - one inside the "Errores" table
# cell ("MSG Expediente no encontrado.") and one in the alternate-trajectory
# paragraph ("MSG Expediente no encontrado."). In both spots we need to turn
# the literal "MSG" into "MSGEX09" (i.e. insert the new error code "EX09"
# right after "MSG", before the following text), while keeping the existing
# "MSG " run's character formatting (Times New Roman, color 002060, etc.) on
# every resulting run.
#
# We look the marker up with Find so the code is independent of absolute
# character offsets, and after each match is handled we resume searching
# from right after it, so every occurrence gets visited exactly once in
# document order.
$rng = $d.Content
while ($rng.Find.Execute("MSG ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $msgStart = $rng.Start
    $fontSize = $rng.Font.Size

    # Insert "EX09" right after "MSG" (before the trailing space / rest of
    # the sentence).
    $msgRng = $d.Range($msgStart, $msgStart + 3)
    $msgRng.InsertAfter("EX09")

    # The inserted text initially shares the exact same run as "MSG" (same
    # formatting), so touch its formatting (toggle Bold on/off, a no-op on
    # the visible result) to force it into its own <w:r> run, split off from
    # "MSG".
    $exRng = $d.Range($msgStart + 3, $msgStart + 7)
    $exRng.Font.Bold = $true
    $exRng.Font.Bold = $false

    if ($fontSize -eq 11) {
        # In the table-cell occurrence, "MSG " was already its own run,
        # separate from the following "Expediente no encontrado." run.
        # Preserve that original run boundary (i.e. keep the single trailing
        # space as its own run) instead of letting it merge back into the
        # following text.
        $spaceRng = $d.Range($msgStart + 7, $msgStart + 8)
        $spaceRng.Font.Bold = $true
        $spaceRng.Font.Bold = $false
    }

    # Continue searching after this occurrence.
    $rng = $d.Range($msgStart + 7, $d.Content.End)
}
